$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries for rows 23 and 24 (Scénáře a Use Case fixes)
$ws.Range("A23").Value = "Oprava 2. iterace - Analytický, BDM, BPM"
$ws.Range("B23").Value = 1

$ws.Range("A24").Value = "Oprava 2. iterace - Úprava scénářů a Use Case"
$ws.Range("B24").Value = 0.5

# Update the active selection to E22 (pane bottomLeft)
$ws.Range("E22").Select()
